$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
  'B2' = 'INSERT INTO ClimateAgreement (CountryCode, CountryName, AgreementName, DateSigned, Target , TargetYear) VALUES (''ARG'', ''Argentina'', ''Paris'', to_date(''2021-11-02 '', ''yyyy-mm-dd''), ''Not exceeding the net emission of 349 million tons of carbon dioxide equivalent in 2030 (target expressed in SAR GWP)'', 2030);'
  'B3' = 'INSERT INTO ClimateAgreement (CountryCode, CountryName, AgreementName, DateSigned, Target , TargetYear) VALUES (''ARG'', ''Argentina'', ''Net Zero'', to_date(''2022-11-06 '', ''yyyy-mm-dd''), ''Argentina plans to make efforts towards reaching GHG emissions neutrality by 2050'', 2050);'
  'B4' = 'INSERT INTO ClimateAgreement (CountryCode, CountryName, AgreementName, DateSigned, Target , TargetYear) VALUES (''AUS'', ''Australia'', ''Paris'', to_date(''2022-06-16 '', ''yyyy-mm-dd''), ''43% below 2005 levels by 2030 (including LULUCF)(Implemented as an absolute target for 2030, and as an emissions budget covering the period 2021-2030)'', 2030);'
  'B5' = 'INSERT INTO ClimateAgreement (CountryCode, CountryName, AgreementName, DateSigned, Target , TargetYear) VALUES (''AUS'', ''Australia'', ''Net Zero'', to_date(''2021-10-29 '', ''yyyy-mm-dd''), ''Australia aims to reach net zero by 2050The target covers all GHG emissions and all sectors of the economy.The target excludes both international aviation and shipping.Australia’s Long-term Emissions Reduction Plan indicates it plans to rely on international and domestic offsets for 10% of reductions  required to reach net zero by 2050.'', 2050);'
  'B6' = 'INSERT INTO ClimateAgreement (CountryCode, CountryName, AgreementName, DateSigned, Target , TargetYear) VALUES (''BTN'', ''Bhutan'', ''Paris'', to_date(''2021-06-01 '', ''yyyy-mm-dd''), ''Remain carbon neutral. International support required to implement additional mitigation measures.'', 2030);'
  'B7' = 'INSERT INTO ClimateAgreement (CountryCode, CountryName, AgreementName, DateSigned, Target , TargetYear) VALUES (''BTN'', ''Bhutan'', ''Net Zero'', to_date(NULL, ''yyyy-mm-dd''), ''Bhutan has a net zero by 2050 goal. At present, the country is already carbon neutral and is committed to remaining that way.'', 2050);'
  'B8' = 'INSERT INTO ClimateAgreement (CountryCode, CountryName, AgreementName, DateSigned, Target , TargetYear) VALUES (''BRA'', ''Brazil'', ''Paris'', to_date(''2023-11-03 '', ''yyyy-mm-dd''), ''2025 net GHG emissions limit of 1.32 GtCO2e (48.4% reduction below 2005 levels)*2030 net GHG emissions limit of 1.20 GtCO2e (53.1% reduction below 2005 levels)'', 2030);'
  'B9' = 'INSERT INTO ClimateAgreement (CountryCode, CountryName, AgreementName, DateSigned, Target , TargetYear) VALUES (''BRA'', ''Brazil'', ''Net Zero'', to_date(NULL, ''yyyy-mm-dd''), ''Climate neutral by 2050'', 2050);'
  'B10' = 'INSERT INTO ClimateAgreement (CountryCode, CountryName, AgreementName, DateSigned, Target , TargetYear) VALUES (''CAN'', ''Canada'', ''Paris'', to_date(''2021-07-12 '', ''yyyy-mm-dd''), ''At least 40-45% below 2005 levels by 2030'', 2030);'
  'B11' = 'INSERT INTO ClimateAgreement (CountryCode, CountryName, AgreementName, DateSigned, Target , TargetYear) VALUES (''CAN'', ''Canada'', ''Net Zero'', to_date(NULL, ''yyyy-mm-dd''), ''Net zero GHG emissions by 2050'', 2050);'
  'B12' = 'INSERT INTO ClimateAgreement (CountryCode, CountryName, AgreementName, DateSigned, Target , TargetYear) VALUES (''CHL'', ''Chile'', ''Paris'', to_date(''2020-04-01 '', ''yyyy-mm-dd''), ''Chile commits to a GHG emission budget not exceeding 1,100 MtCO2e between 2020 and 2030, with a GHG emissions maximum (peak) by 2025, and a GHG emissions level of 95 MtCO2e by 2030.'', 2030);'
  'B13' = 'INSERT INTO ClimateAgreement (CountryCode, CountryName, AgreementName, DateSigned, Target , TargetYear) VALUES (''CHL'', ''Chile'', ''Net Zero'', to_date(''2020-04-01 '', ''yyyy-mm-dd''), ''In addition, under certain specific conditions (financial, markets, technological and political) Chile could exceed a 30% reduction, potentially with a reduction of up to 45% in net emissions by 2030, taking into account actions for GHG emissions mitigation and/or capture.'', 2050);'
  'B14' = 'INSERT INTO ClimateAgreement (CountryCode, CountryName, AgreementName, DateSigned, Target , TargetYear) VALUES (''CHN'', ''China'', ''Paris'', to_date(''2021-10-28 '', ''yyyy-mm-dd''), ''China’s updated NDC contains five overarching targets:1.	Peaking carbon dioxide emissions “before 2030” (up from the previous “around 2030 and making efforts to peak earlier”) and achieve carbon neutrality before 2060.2.	Lower carbon intensity by “over 65%” in 2030 from the 2005 level, (up from the previous “by 60–65%”).3.	Share of non-fossil fuels in primary energy consumption to “around 25%” in 2030, (up from “around 20%”).4.	Increase forest stock volume by around 6 billion cubic metres in 2030 from the 2005 level, (previously 4.5 billion cubic metres).5.	Increase the installed capacity of wind and solar power to over 1,200 GW by 2030 (new target).'', 2030);'
  'B15' = 'INSERT INTO ClimateAgreement (CountryCode, CountryName, AgreementName, DateSigned, Target , TargetYear) VALUES (''CHN'', ''China'', ''Net Zero'', to_date(''2021-10-28 '', ''yyyy-mm-dd''), ''China will strive to reach a CO2emissions peak before 2030, and achieve carbon neutrality before 2060'', 2050);'
  'B16' = 'INSERT INTO ClimateAgreement (CountryCode, CountryName, AgreementName, DateSigned, Target , TargetYear) VALUES (''COL'', ''Colombia'', ''Paris'', to_date(''2020-12-30 '', ''yyyy-mm-dd''), ''Emissions limit of 169.4 MtCO2e in 2030 (equivalent to a 51% reduction below BAU), with an emissions peak by 2027*.Carbon budgets for the period 2020-2030 will be established by 2023.'', 2030);'
  'B17' = 'INSERT INTO ClimateAgreement (CountryCode, CountryName, AgreementName, DateSigned, Target , TargetYear) VALUES (''COL'', ''Colombia'', ''Net Zero'', to_date(''2021-11-21 '', ''yyyy-mm-dd''), ''Carbon neutral by 2050'', 2050);'
  'B18' = 'INSERT INTO ClimateAgreement (CountryCode, CountryName, AgreementName, DateSigned, Target , TargetYear) VALUES (''CRI'', ''Costa-rica'', ''Paris'', to_date(''2020-12-01 '', ''yyyy-mm-dd''), ''Costa Rica commits to an absolute maximum net emissions in 2030 of 9.11 MtCO2e including all emissions and all sectors covered by the corresponding National Greenhouse Gas Emissions Inventory (Government of Costa Rica, 2020).'', 2030);'
  'B19' = 'INSERT INTO ClimateAgreement (CountryCode, CountryName, AgreementName, DateSigned, Target , TargetYear) VALUES (''CRI'', ''Costa-rica'', ''Net Zero'', to_date(''2019-12-01 '', ''yyyy-mm-dd''), ''Costa Rica commits to becoming a decarbonised economy with net-zero emissions by 2050 (Government of Costa Rica, 2019c).'', 2050);'
  'B20' = 'INSERT INTO ClimateAgreement (CountryCode, CountryName, AgreementName, DateSigned, Target , TargetYear) VALUES (''ETH'', ''Ethiopia'', ''Paris'', to_date(''2021-07-23 '', ''yyyy-mm-dd''), ''The unconditional pathway will result in absolute emission levels of 347.3 MtCO2e in 2030 equal to a 14% (-56 MtCO2e) reduction below BAU in 2030. (The 2021 NDC  provides all emissions estimates using global warming potentials of the AR5.)'', 2030);'
  'B21' = 'INSERT INTO ClimateAgreement (CountryCode, CountryName, AgreementName, DateSigned, Target , TargetYear) VALUES (''ETH'', ''Ethiopia'', ''Net Zero'', to_date(''2021-07-23 '', ''yyyy-mm-dd''), ''The conditional pathway will result in absolute emission levels to 125.8 MtCO2e  in 2030 equal to a 68.8% (-277.7 MtCO2e) reduction below BAU in 2030.'', 2050);'
  'B22' = 'INSERT INTO ClimateAgreement (CountryCode, CountryName, AgreementName, DateSigned, Target , TargetYear) VALUES (''DEU'', ''Germany'', ''Paris'', to_date(''2021-08-18 '', ''yyyy-mm-dd''), ''Greenhouse gas emissions will be reduced from 1990 levels as follows:•	at least 65% by 2030'', 2030);'
  'B23' = 'INSERT INTO ClimateAgreement (CountryCode, CountryName, AgreementName, DateSigned, Target , TargetYear) VALUES (''DEU'', ''Germany'', ''Net Zero'', to_date(''2021-08-18 '', ''yyyy-mm-dd''), ''Greenhouse gas emissions will be reduced from 1990 levels as follows:•	at least 88% by 2040By 2045, greenhouse gas emissions will be reduced to such an extent that net greenhouse gas neutrality is achieved. After the year 2050, negative greenhouse gas emissions are to be achieved.'', 2050);'
  'B24' = 'INSERT INTO ClimateAgreement (CountryCode, CountryName, AgreementName, DateSigned, Target , TargetYear) VALUES (''IND'', ''India'', ''Paris'', to_date(''2022-08-26 '', ''yyyy-mm-dd''), ''Emissions intensity of 45% below 2005 levels by 2030'', 2030);'
  'B25' = 'INSERT INTO ClimateAgreement (CountryCode, CountryName, AgreementName, DateSigned, Target , TargetYear) VALUES (''IND'', ''India'', ''Net Zero'', to_date(''2022-08-26 '', ''yyyy-mm-dd''), ''50% cumulative electric power installed capacity from non-fossil fuel-based energy resources by 2030'', 2050);'
  'B26' = 'INSERT INTO ClimateAgreement (CountryCode, CountryName, AgreementName, DateSigned, Target , TargetYear) VALUES (''IDN'', ''Indonesia'', ''Paris'', to_date(''2022-09-23 '', ''yyyy-mm-dd''), ''Reduce emissions by 32% against the 2030 BAU.'', 2030);'
  'B27' = 'INSERT INTO ClimateAgreement (CountryCode, CountryName, AgreementName, DateSigned, Target , TargetYear) VALUES (''IDN'', ''Indonesia'', ''Net Zero'', to_date(''2022-09-23 '', ''yyyy-mm-dd''), ''Reduce emissions by up to 43% against the 2030 BAU.'', 2050);'
  'B28' = 'INSERT INTO ClimateAgreement (CountryCode, CountryName, AgreementName, DateSigned, Target , TargetYear) VALUES (''JPN'', ''Japan'', ''Paris'', to_date(''2021-06-22 '', ''yyyy-mm-dd''), ''46% reduction in 2030 from 2013 levels including LULUCF credits**Japan continues to use a gross-net approach, meaning that Japan does not include the LULUCF sector in its base year (gross) but accounts for net emissions and removals from LULUCF for the target year (net). Japan intends to use LULUCF sink credits up to 47.7 MtCO2e/year.Taking the sector and gas-specific target values in the NDC, this reduces the effectiveness of the 2030 goal from a 46% reduction below 2013 levels including LULUCF to 42% excluding LULUCF. Such an accounting approach undermines the purpose of the Paris Agreement, as it allows for more energy and industry emissions, and should be scrutinised.'', 2030);'
  'B29' = 'INSERT INTO ClimateAgreement (CountryCode, CountryName, AgreementName, DateSigned, Target , TargetYear) VALUES (''JPN'', ''Japan'', ''Net Zero'', to_date(''22/10/2021'', ''yyyy-mm-dd''), ''Carbon neutrality by 2050'', 2050);'
  'B30' = 'INSERT INTO ClimateAgreement (CountryCode, CountryName, AgreementName, DateSigned, Target , TargetYear) VALUES (''KAZ'', ''Kazakhstan'', ''Paris'', to_date(''2023-06-27 '', ''yyyy-mm-dd''), ''Reduction of GHG emissions by 15% by the end of 2030 relative to 1990 base year.'', 2030);'
  'B31' = 'INSERT INTO ClimateAgreement (CountryCode, CountryName, AgreementName, DateSigned, Target , TargetYear) VALUES (''KAZ'', ''Kazakhstan'', ''Net Zero'', to_date(''2023-06-27 '', ''yyyy-mm-dd''), ''Reduction of GHG emissions by 25% by the end of 2030 relative to 1990 base year.'', 2050);'
  'B32' = 'INSERT INTO ClimateAgreement (CountryCode, CountryName, AgreementName, DateSigned, Target , TargetYear) VALUES (''KEN'', ''Kenya'', ''Paris'', to_date(''2020-12-28 '', ''yyyy-mm-dd''), ''Kenya seeks to abate GHG emissions by 32% by 2030 relative to the BAU scenario of 143 MtCO2e and in line with its sustainable development agenda. Kenya intends to bear 21% of the mitigation cost from domestic sources.'', 2030);'
  'B33' = 'INSERT INTO ClimateAgreement (CountryCode, CountryName, AgreementName, DateSigned, Target , TargetYear) VALUES (''KEN'', ''Kenya'', ''Net Zero'', to_date(''2020-12-28 '', ''yyyy-mm-dd''), ''Kenya seeks to abate GHG emissions by 32% by 2030 relative to the BAU scenario of 143 MtCO2e and in line with its sustainable development agenda. 79% of the mitigation cost is subject to international support.'', 2050);'
  'B34' = 'INSERT INTO ClimateAgreement (CountryCode, CountryName, AgreementName, DateSigned, Target , TargetYear) VALUES (''MEX'', ''Mexico'', ''Paris'', to_date(''2022-11-17 '', ''yyyy-mm-dd''), ''Up to 35% GHG below BAU by 2030, where 30% is to be achieved with own resources and the additional 5% with “agreed international support and cooperation for clean energies”.51% reduction of black carbon below BAU by 2030.Note: Here we only consider the 30% percentage to be achieved with own resources.'', 2030);'
  'B35' = 'INSERT INTO ClimateAgreement (CountryCode, CountryName, AgreementName, DateSigned, Target , TargetYear) VALUES (''MEX'', ''Mexico'', ''Net Zero'', to_date(''2022-11-17 '', ''yyyy-mm-dd''), ''Up to 40% GHG, 70% black carbon, below a BAU baseline by 2030.'', 2050);'
  'B36' = 'INSERT INTO ClimateAgreement (CountryCode, CountryName, AgreementName, DateSigned, Target , TargetYear) VALUES (''MAR'', ''Morocco'', ''Paris'', to_date(''2021-06-22 '', ''yyyy-mm-dd''), ''The 2021 NDC increases the unconditional mitigation objective to 18.3% below BAU by 2030'', 2030);'
  'B37' = 'INSERT INTO ClimateAgreement (CountryCode, CountryName, AgreementName, DateSigned, Target , TargetYear) VALUES (''MAR'', ''Morocco'', ''Net Zero'', to_date(''2021-06-22 '', ''yyyy-mm-dd''), ''The 2021 NDC increases the conditional mitigation objective to 45.5% below BAU by 2030'', 2050);'
  'B38' = 'INSERT INTO ClimateAgreement (CountryCode, CountryName, AgreementName, DateSigned, Target , TargetYear) VALUES (''NPL'', ''Nepal'', ''Paris'', to_date(''2020-12-01 '', ''yyyy-mm-dd''), ''Nepal has set an unconditional target of 5000 MW for clean energy generation by 2030, but has no overall target.'', 2030);'
  'B39' = 'INSERT INTO ClimateAgreement (CountryCode, CountryName, AgreementName, DateSigned, Target , TargetYear) VALUES (''NPL'', ''Nepal'', ''Net Zero'', to_date(''2020-12-01 '', ''yyyy-mm-dd''), ''The NDC covers several sectors with target:• Energy sector: expanding clean energy from wind, solar, bioenergy and water 5-10%, from total 15,000 MW.• Transport sector: increasing e-vehicle sales 90% from private passenger for two-wheelers and 60% from all four-wheelers public passenger vehicle.• Residential sector: ensuring 25% using electric stoves for primary cooking and 700 biogas plants (total of household and large scale).• Waste sector:  280 million litters/day of wastewater and 60,000 cubic meters/year of fecal sludge will be managed before released.'', 2050);'
  'B40' = 'INSERT INTO ClimateAgreement (CountryCode, CountryName, AgreementName, DateSigned, Target , TargetYear) VALUES (''NGA'', ''Nigeria'', ''Paris'', to_date(''2021-07-30 '', ''yyyy-mm-dd''), ''20% below BAU by 2030'', 2030);'
  'B41' = 'INSERT INTO ClimateAgreement (CountryCode, CountryName, AgreementName, DateSigned, Target , TargetYear) VALUES (''NGA'', ''Nigeria'', ''Net Zero'', to_date(''2021-07-30 '', ''yyyy-mm-dd''), ''47% below BAU by 2030'', 2050);'
  'B42' = 'INSERT INTO ClimateAgreement (CountryCode, CountryName, AgreementName, DateSigned, Target , TargetYear) VALUES (''NOR'', ''Norway'', ''Paris'', to_date(''2022-11-03 '', ''yyyy-mm-dd''), ''At least 55% reduction in greenhouse gas emissions compared to 1990 levels'', 2030);'
  'B43' = 'INSERT INTO ClimateAgreement (CountryCode, CountryName, AgreementName, DateSigned, Target , TargetYear) VALUES (''NOR'', ''Norway'', ''Net Zero'', to_date(''2020-11-25 '', ''yyyy-mm-dd''), ''Achieve emission reductions of 90–95% from 1990 by 2050'', 2050);'
  'B44' = 'INSERT INTO ClimateAgreement (CountryCode, CountryName, AgreementName, DateSigned, Target , TargetYear) VALUES (''PER'', ''Peru'', ''Paris'', to_date(''2020-12-18 '', ''yyyy-mm-dd''), ''Peru’s national contributions set an unconditional goal of limiting its GHG emissions to a maximum level of 208.8 MtCO2e by 2030'', 2030);'
  'B45' = 'INSERT INTO ClimateAgreement (CountryCode, CountryName, AgreementName, DateSigned, Target , TargetYear) VALUES (''PER'', ''Peru'', ''Net Zero'', to_date(''2020-12-18 '', ''yyyy-mm-dd''), ''Conditioned to the availability of international financing, Peru proposes to limit its GHG emissions to a level of 179 MtCO2e by 2030'', 2050);'
  'B46' = 'INSERT INTO ClimateAgreement (CountryCode, CountryName, AgreementName, DateSigned, Target , TargetYear) VALUES (''PHL'', ''Philippines'', ''Paris'', to_date(''2021-04-15 '', ''yyyy-mm-dd''), ''2.71% of its 75% reduction below a cumulative 2020-2030 BAU trajectory will be achieved unconditionally'', 2030);'
  'B47' = 'INSERT INTO ClimateAgreement (CountryCode, CountryName, AgreementName, DateSigned, Target , TargetYear) VALUES (''PHL'', ''Philippines'', ''Net Zero'', to_date(''2021-04-15 '', ''yyyy-mm-dd''), ''75% below a cumulative 2020-2030 BAU trajectory of which 72.29% is conditional on international support'', 2050);'
  'B48' = 'INSERT INTO ClimateAgreement (CountryCode, CountryName, AgreementName, DateSigned, Target , TargetYear) VALUES (''SGP'', ''Singapore'', ''Paris'', to_date(''2022-11-04 '', ''yyyy-mm-dd''), ''To reduce emissions to around 60 MtCO2e in 2030 after peaking its emissions earlier.'', 2030);'
  'B49' = 'INSERT INTO ClimateAgreement (CountryCode, CountryName, AgreementName, DateSigned, Target , TargetYear) VALUES (''SGP'', ''Singapore'', ''Net Zero'', to_date(''2022-11-03 '', ''yyyy-mm-dd''), ''Net zero GHG emissions by 2050.'', 2050);'
  'B50' = 'INSERT INTO ClimateAgreement (CountryCode, CountryName, AgreementName, DateSigned, Target , TargetYear) VALUES (''CHE'', ''Switzerland'', ''Paris'', to_date(''2021-12-17 '', ''yyyy-mm-dd''), ''A reduction of at least 50% by 2030 compared with 1990 levels, corresponding to an average reduction of greenhouse gas emissions by at least 35% over the period 2021–2030.'', 2030);'
  'B51' = 'INSERT INTO ClimateAgreement (CountryCode, CountryName, AgreementName, DateSigned, Target , TargetYear) VALUES (''CHE'', ''Switzerland'', ''Net Zero'', to_date(''2021-01-28 '', ''yyyy-mm-dd''), ''Switzerland should achieve balanced greenhouse gas performance by 2050 at the latest (net zero).'', 2050);'
  'B52' = 'INSERT INTO ClimateAgreement (CountryCode, CountryName, AgreementName, DateSigned, Target , TargetYear) VALUES (''THA'', ''Thailand'', ''Paris'', to_date(''2022-11-02 '', ''yyyy-mm-dd''), ''Reduce greenhouse gas emissions by 30% from projected business-as-usual level by 2030.'', 2030);'
  'B53' = 'INSERT INTO ClimateAgreement (CountryCode, CountryName, AgreementName, DateSigned, Target , TargetYear) VALUES (''THA'', ''Thailand'', ''Net Zero'', to_date(''2022-11-02 '', ''yyyy-mm-dd''), ''Reduce greenhouse gas emissions by 40% from projected business-as-usual level by 2030, subject to adequate and enhanced access to technology development and transfer, financial resources and capacity building support.'', 2050);'
  'B54' = 'INSERT INTO ClimateAgreement (CountryCode, CountryName, AgreementName, DateSigned, Target , TargetYear) VALUES (''TUR'', ''Turkey'', ''Paris'', to_date(''2023-04-13 '', ''yyyy-mm-dd''), ''41% reduction in GHG emissions from the BAU level by 2030'', 2030);'
  'B55' = 'INSERT INTO ClimateAgreement (CountryCode, CountryName, AgreementName, DateSigned, Target , TargetYear) VALUES (''TUR'', ''Turkey'', ''Net Zero'', to_date(''2021-09-01 '', ''yyyy-mm-dd''), ''Türkiye aims to reach net zero by 2053'', 2050);'
  'B56' = 'INSERT INTO ClimateAgreement (CountryCode, CountryName, AgreementName, DateSigned, Target , TargetYear) VALUES (''VNM'', ''Vietnam'', ''Paris'', to_date(''2022-01-11 '', ''yyyy-mm-dd''), ''15.8% reduction below BAU in 2030'', 2030);'
  'B57' = 'INSERT INTO ClimateAgreement (CountryCode, CountryName, AgreementName, DateSigned, Target , TargetYear) VALUES (''VNM'', ''Vietnam'', ''Net Zero'', to_date(''2022-11-01 '', ''yyyy-mm-dd''), ''43.5% reduction below BAU in 2030'', 2050);'
}

foreach ($cell in $updates.Keys) {
    $ws.Range($cell).Value = $updates[$cell]
}
